$d = $word.ActiveDocument

$replacements = @(
    @{old="377÷9="; new="237÷2="},
    @{old="632÷5="; new="249÷2="},
    @{old="786÷9="; new="679÷3="},
    @{old="206÷8="; new="730÷2="},
    @{old="804÷3="; new="893÷4="},
    @{old="494÷9="; new="252÷5="},
    @{old="921÷8="; new="235÷8="},
    @{old="135÷8="; new="995÷5="},
    @{old="981÷9="; new="965÷9="},
    @{old="205÷8="; new="268÷7="},
    @{old="148÷7="; new="335÷4="},
    @{old="655÷5="; new="264÷6="},
    @{old="907÷7="; new="719÷8="},
    @{old="936÷2="; new="428÷9="},
    @{old="661÷6="; new="601÷2="},
    @{old="931÷7="; new="583÷3="},
    @{old="984÷5="; new="988÷2="},
    @{old="116÷7="; new="243÷9="},
    @{old="841÷2="; new="492÷4="},
    @{old="609÷9="; new="290÷8="},
    @{old="386÷4="; new="594÷2="},
    @{old="791÷9="; new="446÷5="},
    @{old="723÷3="; new="869÷8="},
    @{old="654÷6="; new="159÷9="},
    @{old="947÷6="; new="948÷5="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $r.new, 2)
}
